# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45203 (2023-10-04) to 45204 (2023-10-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows are 2 through 530 (row 1 is the header row).
$ws.Range("C2:C530").Value = 45204
